$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.248.22'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '3.600.39'
$ws.Range('E3').Value = '  +1.73%  '
$ws.Range('D5').Value = '603.55'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = '195.89'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('D11').Value = '53.77'
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').Value = '4.174.91'
$ws.Range('E14').Value = '  +2.08%  '
$ws.Range('D15').Value = '13.07'
$ws.Range('E15').Value = '  +3.06%  '
$ws.Range('D16').Value = '595.94'
$ws.Range('E16').Value = '  -0.66%  '
$ws.Range('D17').Value = '70.346.68'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = '3.611.41'
$ws.Range('E18').Value = '  +2.44%  '
$ws.Range('D19').Value = '19.05'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('E23').Value = '  -2.08%  '
$ws.Range('D24').Value = '102.03'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('D28').Value = '9.64'
$ws.Range('E28').Value = '  -0.41%  '
$ws.Range('D29').Value = '33.77'
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('D30').Value = '4.79'
$ws.Range('E30').Value = '  +6.94%  '
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('E32').Value = '  -3.78%  '
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('B34').Value = 'PEPE'
$ws.Range('C34').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D34').Value = '0.0₃0903'
$ws.Range('E34').Value = '  +8.92%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '63.19'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = '3.898.72'
$ws.Range('E36').Value = '  +4.15%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').Value = '522.12'
$ws.Range('E39').Value = '  +4.78%  '
$ws.Range('D40').Value = '36.91'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('E42').Value = '  -2.36%  '
$ws.Range('E43').Value = '  -2.08%  '
$ws.Range('D44').Value = '0.0453'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D45').Value = '3.41'
$ws.Range('E45').Value = '  +2.73%  '
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  -0.67%  '
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('D50').Value = '0.000252'
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('E51').Value = '  +0.26%  '
